$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(42, 8).Value = 409.85715
$ws.Cells.Item(42, 9).Value = 515
$ws.Cells.Item(42, 10).Value = 367.8
$ws.Cells.Item(42, 11).Value = 1545
$ws.Cells.Item(42, 12).Value = 1103.4
$ws.Cells.Item(42, 13).Value = -1315
$ws.Cells.Item(42, 14).Value = -1563.4
$ws.Cells.Item(98, 8).Value = 10616.923
$ws.Cells.Item(98, 9).Value = 5602
$ws.Cells.Item(98, 10).Value = 27333.334
$ws.Cells.Item(98, 11).Value = 5602
$ws.Cells.Item(98, 12).Value = 27333.334
$ws.Cells.Item(98, 13).Value = -4104
$ws.Cells.Item(98, 14).Value = -30329.334
$ws.Cells.Item(113, 8).Value = 3089.6428
$ws.Cells.Item(113, 9).Value = 2463.75
$ws.Cells.Item(113, 10).Value = 3340
$ws.Cells.Item(113, 11).Value = 2463.75
$ws.Cells.Item(113, 12).Value = 3340
$ws.Cells.Item(113, 13).Value = 790.25
$ws.Cells.Item(113, 14).Value = -9848
$ws.Cells.Item(122, 8).Value = 10616.923
$ws.Cells.Item(122, 9).Value = 5602
$ws.Cells.Item(122, 10).Value = 27333.334
$ws.Cells.Item(122, 11).Value = 16806
$ws.Cells.Item(122, 12).Value = 82000.00199999999
$ws.Cells.Item(122, 13).Value = -14356
$ws.Cells.Item(122, 14).Value = -86900.00199999999
$ws.Cells.Item(123, 8).Value = 67098.17999999999
$ws.Cells.Item(123, 10).Value = 67098.17999999999
$ws.Cells.Item(123, 12).Value = 67098.17999999999
$ws.Cells.Item(123, 14).Value = -76898.17999999999
$ws.Cells.Item(124, 8).Value = 0
$ws.Cells.Item(124, 10).Value = 0
$ws.Cells.Item(124, 12).Value = 0
$ws.Cells.Item(126, 8).Value = 46796.668
$ws.Cells.Item(126, 10).Value = 46796.668
$ws.Cells.Item(126, 12).Value = 46796.668
$ws.Cells.Item(126, 14).Value = -56676.668
$ws.Cells.Item(132, 8).Value = 3179.4
$ws.Cells.Item(132, 9).Value = 3178.3157
$ws.Cells.Item(132, 11).Value = 9534.947100000001
$ws.Cells.Item(132, 13).Value = -7004.947100000001
$ws.Cells.Item(137, 8).Value = 1714.3077
$ws.Cells.Item(137, 10).Value = 1750
$ws.Cells.Item(137, 12).Value = 5250
$ws.Cells.Item(137, 14).Value = -10350
$ws.Cells.Item(141, 8).Value = 4997.9653
$ws.Cells.Item(141, 9).Value = 1605.0358
$ws.Cells.Item(141, 11).Value = 4815.107400000001
$ws.Cells.Item(141, 13).Value = 364.8925999999992
$ws.Cells.Item(124, 14).ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 15845.775
$ws.Cells.Item(32, 9).Value = 17483.594
$ws.Cells.Item(32, 10).Value = 9294.5
$ws.Cells.Item(32, 11).Value = 17483.594
$ws.Cells.Item(32, 12).Value = 9294.5
$ws.Cells.Item(32, 13).Value = -17196.594
$ws.Cells.Item(32, 14).Value = -9868.5
$ws.Cells.Item(122, 8).Value = 2504.3076
$ws.Cells.Item(122, 9).Value = 2033.55
$ws.Cells.Item(122, 11).Value = 6100.65
$ws.Cells.Item(122, 13).Value = -3650.65

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 127790
$ws.Cells.Item(86, 9).Value = 2200
$ws.Cells.Item(86, 10).Value = 169653.33
$ws.Cells.Item(86, 11).Value = 2200
$ws.Cells.Item(86, 12).Value = 169653.33
$ws.Cells.Item(86, 13).Value = -1077
$ws.Cells.Item(86, 14).Value = -171899.33
$ws.Cells.Item(89, 8).Value = 127790
$ws.Cells.Item(89, 9).Value = 2200
$ws.Cells.Item(89, 10).Value = 169653.33
$ws.Cells.Item(89, 11).Value = 11000
$ws.Cells.Item(89, 12).Value = 848266.6499999999
$ws.Cells.Item(89, 13).Value = -5384
$ws.Cells.Item(89, 14).Value = -859498.6499999999
$ws.Cells.Item(134, 8).Value = 2157.3103
$ws.Cells.Item(134, 9).Value = 1813.25
$ws.Cells.Item(134, 11).Value = 5439.75
$ws.Cells.Item(134, 13).Value = -2904.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3263.25
$ws.Cells.Item(31, 9).Value = 1958.8572
$ws.Cells.Item(31, 10).Value = 4277.778
$ws.Cells.Item(31, 11).Value = 1958.8572
$ws.Cells.Item(31, 12).Value = 4277.778
$ws.Cells.Item(31, 13).Value = -1663.8572
$ws.Cells.Item(31, 14).Value = -4867.778
$ws.Cells.Item(34, 8).Value = 3263.25
$ws.Cells.Item(34, 9).Value = 1958.8572
$ws.Cells.Item(34, 10).Value = 4277.778
$ws.Cells.Item(34, 11).Value = 1958.8572
$ws.Cells.Item(34, 12).Value = 4277.778
$ws.Cells.Item(34, 13).Value = -1756.8572
$ws.Cells.Item(34, 14).Value = -4681.778
$ws.Cells.Item(62, 8).Value = 85867.5
$ws.Cells.Item(62, 9).Value = 127026.25
$ws.Cells.Item(62, 10).Value = 3550
$ws.Cells.Item(62, 11).Value = 127026.25
$ws.Cells.Item(62, 12).Value = 3550
$ws.Cells.Item(62, 13).Value = -126402.25
$ws.Cells.Item(62, 14).Value = -4798
$ws.Cells.Item(65, 8).Value = 85867.5
$ws.Cells.Item(65, 9).Value = 127026.25
$ws.Cells.Item(65, 10).Value = 3550
$ws.Cells.Item(65, 11).Value = 635131.25
$ws.Cells.Item(65, 12).Value = 17750
$ws.Cells.Item(65, 13).Value = -632011.25
$ws.Cells.Item(65, 14).Value = -23990
$ws.Cells.Item(122, 8).Value = 1972.2142
$ws.Cells.Item(122, 9).Value = 1933.7778
$ws.Cells.Item(122, 10).Value = 2041.4
$ws.Cells.Item(122, 11).Value = 5801.3334
$ws.Cells.Item(122, 12).Value = 6124.200000000001
$ws.Cells.Item(122, 13).Value = -3351.3334
$ws.Cells.Item(122, 14).Value = -11024.2
$ws.Cells.Item(134, 8).Value = 2550.0952
$ws.Cells.Item(134, 9).Value = 1824.9231
$ws.Cells.Item(134, 10).Value = 3728.5
$ws.Cells.Item(134, 11).Value = 5474.7693
$ws.Cells.Item(134, 12).Value = 11185.5
$ws.Cells.Item(134, 13).Value = -2939.7693
$ws.Cells.Item(134, 14).Value = -16255.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 291.38095
$ws.Cells.Item(107, 9).Value = 347.91306
$ws.Cells.Item(107, 10).Value = 222.94737
$ws.Cells.Item(107, 11).Value = 1043.73918
$ws.Cells.Item(107, 12).Value = 668.84211
$ws.Cells.Item(107, 13).Value = 876.26082
$ws.Cells.Item(107, 14).Value = -4508.84211
$ws.Cells.Item(132, 8).Value = 2289.1875
$ws.Cells.Item(132, 9).Value = 1853
$ws.Cells.Item(132, 11).Value = 16677
$ws.Cells.Item(132, 13).Value = -14147
$ws.Cells.Item(133, 8).Value = 5807.273
$ws.Cells.Item(133, 9).Value = 1990.5
$ws.Cells.Item(133, 10).Value = 6655.4443
$ws.Cells.Item(133, 11).Value = 5971.5
$ws.Cells.Item(133, 12).Value = 19966.3329
$ws.Cells.Item(133, 13).Value = -911.5
$ws.Cells.Item(133, 14).Value = -30086.3329
$ws.Cells.Item(141, 8).Value = 5838
$ws.Cells.Item(141, 9).Value = 5838
$ws.Cells.Item(141, 11).Value = 17514
$ws.Cells.Item(141, 13).Value = -12334

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 3631
$ws.Cells.Item(122, 9).Value = 2552.3333
$ws.Cells.Item(122, 10).Value = 3954.6
$ws.Cells.Item(122, 11).Value = 7656.999899999999
$ws.Cells.Item(122, 12).Value = 11863.8
$ws.Cells.Item(122, 13).Value = -5206.999899999999
$ws.Cells.Item(122, 14).Value = -16763.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5020
$ws.Cells.Item(7, 9).Value = 5275
$ws.Cells.Item(7, 11).Value = 5275
$ws.Cells.Item(7, 13).Value = -5163
$ws.Cells.Item(40, 8).Value = 4006.1667
$ws.Cells.Item(40, 9).Value = 5379
$ws.Cells.Item(40, 10).Value = 2633.3333
$ws.Cells.Item(40, 11).Value = 5379
$ws.Cells.Item(40, 12).Value = 2633.3333
$ws.Cells.Item(40, 13).Value = -5243
$ws.Cells.Item(40, 14).Value = -2905.3333
$ws.Cells.Item(61, 8).Value = 3233.3333
$ws.Cells.Item(61, 9).Value = 7000
$ws.Cells.Item(61, 10).Value = 1350
$ws.Cells.Item(61, 11).Value = 7000
$ws.Cells.Item(61, 12).Value = 1350
$ws.Cells.Item(61, 13).Value = -6798
$ws.Cells.Item(61, 14).Value = -1754
$ws.Cells.Item(82, 8).Value = 4306.25
$ws.Cells.Item(82, 10).Value = 5075
$ws.Cells.Item(82, 12).Value = 5075
$ws.Cells.Item(82, 14).Value = -5797
$ws.Cells.Item(85, 8).Value = 4306.25
$ws.Cells.Item(85, 10).Value = 5075
$ws.Cells.Item(85, 12).Value = 5075
$ws.Cells.Item(85, 14).Value = -7571
$ws.Cells.Item(94, 8).Value = 10564
$ws.Cells.Item(94, 10).Value = 10564
$ws.Cells.Item(94, 12).Value = 10564
$ws.Cells.Item(94, 14).Value = -11916
$ws.Cells.Item(113, 8).Value = 3233.3333
$ws.Cells.Item(113, 9).Value = 7000
$ws.Cells.Item(113, 10).Value = 1350
$ws.Cells.Item(113, 11).Value = 7000
$ws.Cells.Item(113, 12).Value = 1350
$ws.Cells.Item(113, 13).Value = -4830
$ws.Cells.Item(113, 14).Value = -5690
$ws.Cells.Item(121, 8).Value = 63856.93
$ws.Cells.Item(121, 10).Value = 63856.93
$ws.Cells.Item(121, 12).Value = 63856.93
$ws.Cells.Item(121, 14).Value = -67350.92999999999
$ws.Cells.Item(122, 8).Value = 22227000
$ws.Cells.Item(122, 9).Value = 4714.857
$ws.Cells.Item(122, 10).Value = 100005000
$ws.Cells.Item(122, 11).Value = 14144.571
$ws.Cells.Item(122, 12).Value = 300015000
$ws.Cells.Item(122, 13).Value = -11694.571
$ws.Cells.Item(122, 14).Value = -300019900
$ws.Cells.Item(126, 8).Value = 5020
$ws.Cells.Item(126, 9).Value = 5275
$ws.Cells.Item(126, 11).Value = 15825
$ws.Cells.Item(126, 13).Value = -13355

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 563.0769
$ws.Cells.Item(113, 9).Value = 611.8182
$ws.Cells.Item(113, 10).Value = 295
$ws.Cells.Item(113, 11).Value = 1835.4546
$ws.Cells.Item(113, 12).Value = 885
$ws.Cells.Item(113, 13).Value = 334.5454
$ws.Cells.Item(113, 14).Value = -5225
$ws.Cells.Item(122, 8).Value = 31252084
$ws.Cells.Item(122, 9).Value = 47620420
$ws.Cells.Item(122, 11).Value = 142861260
$ws.Cells.Item(122, 13).Value = -142858810
$ws.Cells.Item(126, 8).Value = 5658.4116
$ws.Cells.Item(126, 9).Value = 5949.75
$ws.Cells.Item(126, 11).Value = 17849.25
$ws.Cells.Item(126, 13).Value = -15379.25
